$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the data block so that numeric-looking strings
# (e.g. "1", "0,1960") are stored as shared-string text, not as numbers,
# matching the authored workbook which keeps these as text cells (t="s").
$ws.Range("A2:I9").NumberFormat = "@"

$ws.Range("A2").Value = "    BRL"
$ws.Range("B2").Value = "1"
$ws.Range("C2").Value = "0,1960"
$ws.Range("D2").Value = "0,1761"
$ws.Range("E2").Value = "0,1574"
$ws.Range("F2").Value = "21,1615"
$ws.Range("G2").Value = "0,1866"
$ws.Range("H2").Value = "0,2741"
$ws.Range("I2").Value = "0,3180"
$ws.Range("A3").Value = "    USD"
$ws.Range("B3").Value = "5,1010"
$ws.Range("C3").Value = "1"
$ws.Range("D3").Value = "0,8974"
$ws.Range("E3").Value = "0,8040"
$ws.Range("F3").Value = "107,94"
$ws.Range("G3").Value = "0,9513"
$ws.Range("H3").Value = "1,3979"
$ws.Range("I3").Value = "1,6221"
$ws.Range("A4").Value = "    EUR"
$ws.Range("B4").Value = "5,6836"
$ws.Range("C4").Value = "1,1143"
$ws.Range("D4").Value = "1"
$ws.Range("E4").Value = "0,8958"
$ws.Range("F4").Value = "120,27"
$ws.Range("G4").Value = "1,0604"
$ws.Range("H4").Value = "1,5581"
$ws.Range("I4").Value = "1,8082"
$ws.Range("A5").Value = "    GBP"
$ws.Range("B5").Value = "6,3548"
$ws.Range("C5").Value = "1,2439"
$ws.Range("D5").Value = "1,1180"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "134,26"
$ws.Range("G5").Value = "1,1834"
$ws.Range("H5").Value = "1,7394"
$ws.Range("I5").Value = "2,0177"
$ws.Range("A6").Value = "    JPY"
$ws.Range("B6").Value = "0,04725"
$ws.Range("C6").Value = "0,00926"
$ws.Range("D6").Value = "0,00832"
$ws.Range("E6").Value = "0,00744"
$ws.Range("F6").Value = "1"
$ws.Range("G6").Value = "0,0088"
$ws.Range("H6").Value = "0,01295"
$ws.Range("I6").Value = "0,01504"
$ws.Range("A7").Value = "    CHF"
$ws.Range("B7").Value = "5,3602"
$ws.Range("C7").Value = "1,0512"
$ws.Range("D7").Value = "0,9431"
$ws.Range("E7").Value = "0,8449"
$ws.Range("F7").Value = "113,47"
$ws.Range("G7").Value = "1"
$ws.Range("H7").Value = "1,4700"
$ws.Range("I7").Value = "1,7051"
$ws.Range("A8").Value = "    CAD"
$ws.Range("B8").Value = "3,6479"
$ws.Range("C8").Value = "0,7154"
$ws.Range("D8").Value = "0,6420"
$ws.Range("E8").Value = "0,5743"
$ws.Range("F8").Value = "77,19"
$ws.Range("G8").Value = "0,6806"
$ws.Range("H8").Value = "1"
$ws.Range("I8").Value = "1,1602"
$ws.Range("A9").Value = "    AUD"
$ws.Range("B9").Value = "3,1449"
$ws.Range("C9").Value = "0,6165"
$ws.Range("D9").Value = "0,5534"
$ws.Range("E9").Value = "0,4955"
$ws.Range("F9").Value = "66,55"
$ws.Range("G9").Value = "0,5865"
$ws.Range("H9").Value = "0,8620"
$ws.Range("I9").Value = "1"

# Remove the style index (s="1") that was inherited from the old blank
# row 2 cells so the new data rows have no explicit cell style, matching
# the target sheet XML.
$ws.Range("A2:I9").ClearFormats()
